$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DropTable")

# Insert a new row above current row 3 (the 5001 entry), pushing it down to row 4
$ws.Rows("3").Insert()

# Populate the new row 3 with the new drop entry (dropId 1002, enum Heart)
$ws.Range("A3").Value = 1002
$ws.Range("B3").Value = "Heart"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2

# Update the dropEnum text for the (now) row 4 entry to the capitalized form
$ws.Range("B4").Value = "Exp, Gold, Heart, Gacha"
